$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat'
$ws.Range('G3').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud'
$ws.Range('G4').Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G5').Value = 'Dr. Eman Tantawi, Administrator, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid'
$ws.Range('G6').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G8').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range('G9').Value = 'Dr. Nourhan Mohammad, Dr. Safa Hany, D Wessam Atef, Dr. Sara Nabil'
$ws.Range('G10').Value = 'Dr. Sara Nabil, Dr. Aya Saeed, D Wessam Atef, Dr. Omnia Mohammad, Dr. Amal Awwad'
$ws.Range('G12').Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Maryam Ahmad, Dr. Merna Mahrous'
$ws.Range('G13').Value = 'Dr. Eman M. Elsaid, Dr. Dina Adel, Dr. Nourhan Osama, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim'
$ws.Range('G14').Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range('G15').Value = 'Dr. Mayar Ahmad Embaby, Nourhan Mamdouh Hassan, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range('G16').Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Al-Shimaa Khaled'
$ws.Range('G23').Value = 'Menna tuâ€™Allah Gamil, Dr. Mona Ibrahim Hussein'
$ws.Range('G25').Value = 'Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein'
$ws.Range('G27').Value = 'Administrator, Dr. Rada Rabea, Dr. Marian Samir, Dr. Hana Amr, Dr. Nourhan Mohammad'
$ws.Range('G29').Value = 'Dr. Marina Atef, Dr. Remon, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nahla, Dr. Marina Sorial, Dr. Nardine, Dr. Aya Hanafy, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range('G30').Value = 'Dr. Nardine, Dr. Nahla'
$ws.Range('G33').Value = 'Dr. Servinaz Sayed Mohammad, Dr. Menna tuâ€™Allah Medhat'
$ws.Range('G34').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud'
$ws.Range('G35').Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range('G36').Value = 'Dr. Eman Tantawi, Administrator, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid'
$ws.Range('G37').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G38').Value = 'Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range('G40').Value = 'Dr. Nourhan Mohammad, Dr. Safa Hany, D Wessam Atef, Dr. Sara Nabil'
$ws.Range('G41').Value = 'Dr. Amal Awwad, Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range('G43').Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Maryam Ahmad, Dr. Merna Mahrous'
$ws.Range('G44').Value = 'Dr. Eman M. Elsaid, Dr. Dina Adel, Dr. Nourhan Osama, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim'
$ws.Range('G45').Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range('G52').Value = 'Administrator, Dr. Afnan Fares'
$ws.Range('G55').Value = 'Menna tuâ€™Allah Gamil, Dr. Mona Ibrahim Hussein'
$ws.Range('G56').Value = 'Dr. Youstina Gamil, Dr. Mona Ibrahim Hussein'
$ws.Range('G58').Value = 'Administrator, Dr. Rada Rabea, Dr. Marian Samir, Dr. Hana Amr, Dr. Nourhan Mohammad'
$ws.Range('G60').Value = 'Dr. Marina Atef, Dr. Remon, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nahla, Dr. Marina Sorial, Dr. Nardine, Dr. Aya Hanafy, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range('G61').Value = 'Dr. Nardine, Dr. Nahla'
$ws.Range('G64').Value = 'Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat'
$ws.Range('G65').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G66').Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Hanan Ragab'
$ws.Range('G67').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G68').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G69').Value = 'Dr. Fatma Elhady, Dr. Kerelos Zareef'
$ws.Range('G70').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range('G71').Value = 'Dr. Mariam Nour El-Din, Dr. Safa Hany, Dr. Omnia Mohammad, Dr. Sara Nabil'
$ws.Range('G72').Value = 'Dr. Sara Nabil, Dr. Aya Saeed, D Wessam Atef, Dr. Omnia Mohammad, Dr. Amal Awwad'
$ws.Range('G74').Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Merna Mahrous'
$ws.Range('G75').Value = 'Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Amany Raafat'
$ws.Range('G76').Value = 'Dr. Mayar Ahmad Embaby, Nourhan Mamdouh Hassan, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range('G89').Value = 'Dr. Nourham Mostafa, Dr. Aya Alaa-Eldein'
$ws.Range('G91').Value = 'Dr. Marina Atef, Dr. Remon, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nahla, Dr. Marina Sorial, Dr. Nardine, Dr. Aya Hanafy, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range('G92').Value = 'Dr. Nardine, Dr. Nahla'
$ws.Range('G95').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range('G96').Value = 'Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud'
$ws.Range('G97').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range('G98').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range('G99').Value = 'Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub'
$ws.Range('G101').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range('G102').Value = 'Dr. Mariam Nour El-Din, Dr. Nourhan Mohammad, Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G103').Value = 'Dr. Amal Awwad, Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range('G105').Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Maryam Ahmad, Dr. Merna Mahrous'
$ws.Range('G106').Value = 'Dr. Dina Adel, Dr. Nadia Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh'
$ws.Range('G107').Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range('G108').Value = 'Dr. Mayar Ahmad Embaby, Nourhan Mamdouh Hassan, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range('G109').Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Al-Shimaa Khaled'
$ws.Range('G114').Value = 'Administrator, Dr. Afnan Fares'
$ws.Range('G116').Value = 'Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G117').Value = 'Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G120').Value = 'Dr. Ahmad Mostafa, Dr. Marian Samir, Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Nourhan Mohammad'
$ws.Range('G122').Value = 'Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range('G123').Value = 'Dr. Nardine, Dr. Nahla'
$ws.Range('G126').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range('G127').Value = 'Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Nourhan Mahmoud'
$ws.Range('G128').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range('G129').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud, Dr. Nahla Nagiub, Dr. Asmaa Reda'
$ws.Range('G130').Value = 'Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Nahla Nagiub'
$ws.Range('G131').Value = 'Dr. Fatma Elhady, Dr. Kerelos Zareef'
$ws.Range('G132').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range('G133').Value = 'Dr. Mariam Nour El-Din, Dr. Nourhan Mohammad, Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G134').Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range('G136').Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Maryam Ahmad, Dr. Merna Mahrous'
$ws.Range('G137').Value = 'Dr. Eman M. Elsaid, Dr. Dina Adel, Dr. Nourhan Osama, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim'
$ws.Range('G138').Value = 'Dr. Mayar Ahmad Embaby, Nourhan Mamdouh Hassan, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range('G139').Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range('G147').Value = 'Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G148').Value = 'Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range('G151').Value = 'Dr. Taqwa Mohammad, Dr. Amr Saeed, Dr. Enas Omran'
$ws.Range('G153').Value = 'Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range('G154').Value = 'Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Youstina Magdy'
$ws.Range('G157').Value = 'Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Hend Mahmoud'
$ws.Range('G158').Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G159').Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range('G160').Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G161').Value = 'Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range('G162').Value = 'Dr. Fatma Elhady, Dr. Kerelos Zareef'
$ws.Range('G164').Value = 'Dr. Mariam Nour El-Din, Dr. Safa Hany, Dr. Omnia Mohammad, Dr. Sara Nabil'
$ws.Range('G165').Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range('G167').Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Merna Mahrous'
$ws.Range('G168').Value = 'Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Amany Raafat'
$ws.Range('G169').Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range('G170').Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range('G179').Value = 'Menna tuâ€™Allah Gamil, Dr. Mona Ibrahim Hussein'
$ws.Range('G182').Value = 'Dr. Ahmad Mostafa, Dr. Marian Samir, Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Nourhan Mohammad'
$ws.Range('G184').Value = 'Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nahla, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida'
$ws.Range('G185').Value = 'Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Youstina Magdy'
